$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.62%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.77%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.095'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.08%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05590'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.25%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.466'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.22%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8135'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.12%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8441'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.06%'
$ws.Range("B9").Value = 'MandalaExchangeToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06952'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.11%'
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.02817'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.59%'
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09387'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.13%'
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.001514'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.07%'
$ws.Range("B13").Value = 'TigerCash'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.006215'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.36%'
$ws.Range("B14").Value = 'LEO'
$ws.Range("C14").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.607'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.05%'
$ws.Range("B15").Value = 'GateToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.021'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.12%'
$ws.Range("B16").Value = 'BTSEToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.055'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.74%'
$ws.Range("B17").Value = 'BitpandaEcosystemToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.3112'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.33%'
$ws.Range("B18").Value = 'WazirX'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.1335'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.12%'
$ws.Range("B19").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C19").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.03183'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.64%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1294'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.96%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.762'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.28%'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04649'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.41%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1376'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.57%'
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0005996'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.60%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001246'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.23%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004551'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '6.16%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009610'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.92%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001940'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.00%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03651'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.24%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006187'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '83.41%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-22.43%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002590'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.09%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008068'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.55%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005386'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.69%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.10%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1452'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-19.36%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002422'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '20.13%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002102'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.10%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002002'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.10%'
